$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Połączenia")

# Update the constraint-equation text and move it from H5 to H2.
$ws.Range("H2").Value = "t^2+pi/2"
$ws.Range("H5").ClearContents()

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("H18").Select()
